$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the "From" value of rule R30 (cell C10) from 18 to 1
$ws.Range("C10").Value = 1
